$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 432-463: shift/rotate of weekly price-quality rows (Especial/Primera/Segunda)
# columns changed: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$ws.Range("D432").Value = 44826
$ws.Range("L432").Value = "Primera"
$ws.Range("M432").Value = 480
$ws.Range("N432").Value = 19000
$ws.Range("O432").Value = 20000
$ws.Range("P432").Value = 19500
$ws.Range("S432").Value = 2786

$ws.Range("D433").Value = 44826
$ws.Range("L433").Value = "Segunda"
$ws.Range("M433").Value = 320
$ws.Range("N433").Value = 14000
$ws.Range("O433").Value = 15000
$ws.Range("P433").Value = 14500
$ws.Range("S433").Value = 2071

$ws.Range("D434").Value = 44328
$ws.Range("L434").Value = "Especial"
$ws.Range("M434").Value = 200
$ws.Range("N434").Value = 19500
$ws.Range("O434").Value = 20000
$ws.Range("P434").Value = 19750
$ws.Range("S434").Value = 2821

$ws.Range("D435").Value = 44328
$ws.Range("L435").Value = "Primera"
$ws.Range("M435").Value = 300
$ws.Range("N435").Value = 16500
$ws.Range("O435").Value = 17000
$ws.Range("P435").Value = 16750
$ws.Range("S435").Value = 2393

$ws.Range("D436").Value = 44328
$ws.Range("L436").Value = "Segunda"
$ws.Range("M436").Value = 200
$ws.Range("N436").Value = 12500
$ws.Range("O436").Value = 13000
$ws.Range("P436").Value = 12750
$ws.Range("S436").Value = 1821

$ws.Range("D437").Value = 44356
$ws.Range("L437").Value = "Especial"
$ws.Range("M437").Value = 200
$ws.Range("N437").Value = 17500
$ws.Range("O437").Value = 18000
$ws.Range("P437").Value = 17750
$ws.Range("S437").Value = 2536

$ws.Range("D438").Value = 44356
$ws.Range("L438").Value = "Primera"
$ws.Range("M438").Value = 300
$ws.Range("N438").Value = 15500
$ws.Range("O438").Value = 16000
$ws.Range("P438").Value = 15750
$ws.Range("S438").Value = 2250

$ws.Range("D439").Value = 44356
$ws.Range("L439").Value = "Segunda"
$ws.Range("M439").Value = 200
$ws.Range("N439").Value = 12500
$ws.Range("O439").Value = 13000
$ws.Range("P439").Value = 12750
$ws.Range("S439").Value = 1821

$ws.Range("D440").Value = 44322
$ws.Range("L440").Value = "Especial"
$ws.Range("M440").Value = 300
$ws.Range("N440").Value = 16500
$ws.Range("O440").Value = 17000
$ws.Range("P440").Value = 16750
$ws.Range("S440").Value = 2393

$ws.Range("D441").Value = 44322
$ws.Range("L441").Value = "Primera"
$ws.Range("M441").Value = 300
$ws.Range("N441").Value = 14500
$ws.Range("O441").Value = 15000
$ws.Range("P441").Value = 14750
$ws.Range("S441").Value = 2107

$ws.Range("D442").Value = 44322
$ws.Range("L442").Value = "Segunda"
$ws.Range("M442").Value = 240
$ws.Range("N442").Value = 11500
$ws.Range("O442").Value = 12000
$ws.Range("P442").Value = 11750
$ws.Range("S442").Value = 1679

$ws.Range("D443").Value = 44497
$ws.Range("L443").Value = "Especial"
$ws.Range("M443").Value = 300
$ws.Range("N443").Value = 12500
$ws.Range("O443").Value = 13000
$ws.Range("P443").Value = 12750
$ws.Range("S443").Value = 1821

$ws.Range("D444").Value = 44497
$ws.Range("L444").Value = "Primera"
$ws.Range("M444").Value = 400
$ws.Range("N444").Value = 10500
$ws.Range("O444").Value = 11000
$ws.Range("P444").Value = 10750
$ws.Range("S444").Value = 1536

$ws.Range("D445").Value = 44497
$ws.Range("L445").Value = "Segunda"
$ws.Range("M445").Value = 400
$ws.Range("N445").Value = 8500
$ws.Range("O445").Value = 9000
$ws.Range("P445").Value = 8750
$ws.Range("S445").Value = 1250

$ws.Range("D446").Value = 44435
$ws.Range("L446").Value = "Primera"
$ws.Range("M446").Value = 540
$ws.Range("N446").Value = 27500
$ws.Range("O446").Value = 28000
$ws.Range("P446").Value = 27750
$ws.Range("S446").Value = 3964

$ws.Range("D447").Value = 44435
$ws.Range("L447").Value = "Segunda"
$ws.Range("M447").Value = 480
$ws.Range("N447").Value = 21000
$ws.Range("O447").Value = 22000
$ws.Range("P447").Value = 21500
$ws.Range("S447").Value = 3071

$ws.Range("D448").Value = 44251
$ws.Range("L448").Value = "Especial"
$ws.Range("M448").Value = 240
$ws.Range("N448").Value = 15500
$ws.Range("O448").Value = 16000
$ws.Range("P448").Value = 15750
$ws.Range("S448").Value = 2250

$ws.Range("D449").Value = 44251
$ws.Range("L449").Value = "Primera"
$ws.Range("M449").Value = 300
$ws.Range("N449").Value = 13500
$ws.Range("O449").Value = 14000
$ws.Range("P449").Value = 13750
$ws.Range("S449").Value = 1964

$ws.Range("D450").Value = 44251
$ws.Range("L450").Value = "Segunda"
$ws.Range("M450").Value = 300
$ws.Range("N450").Value = 10500
$ws.Range("O450").Value = 11000
$ws.Range("P450").Value = 10750
$ws.Range("S450").Value = 1536

$ws.Range("D451").Value = 44455
$ws.Range("L451").Value = "Especial"
$ws.Range("M451").Value = 200
$ws.Range("N451").Value = 29000
$ws.Range("O451").Value = 30000
$ws.Range("P451").Value = 29500
$ws.Range("S451").Value = 4214

$ws.Range("D452").Value = 44455
$ws.Range("L452").Value = "Primera"
$ws.Range("M452").Value = 300
$ws.Range("N452").Value = 24000
$ws.Range("O452").Value = 25000
$ws.Range("P452").Value = 24500
$ws.Range("S452").Value = 3500

$ws.Range("D453").Value = 44455
$ws.Range("L453").Value = "Segunda"
$ws.Range("M453").Value = 240
$ws.Range("N453").Value = 19000
$ws.Range("O453").Value = 20000
$ws.Range("P453").Value = 19500
$ws.Range("S453").Value = 2786

$ws.Range("D454").Value = 44504
$ws.Range("L454").Value = "Especial"
$ws.Range("M454").Value = 400
$ws.Range("N454").Value = 12500
$ws.Range("O454").Value = 13000
$ws.Range("P454").Value = 12750
$ws.Range("S454").Value = 1821

$ws.Range("D455").Value = 44504
$ws.Range("L455").Value = "Primera"
$ws.Range("M455").Value = 340
$ws.Range("N455").Value = 10500
$ws.Range("O455").Value = 11000
$ws.Range("P455").Value = 10750
$ws.Range("S455").Value = 1536

$ws.Range("D456").Value = 44504
$ws.Range("L456").Value = "Segunda"
$ws.Range("M456").Value = 240
$ws.Range("N456").Value = 8500
$ws.Range("O456").Value = 9000
$ws.Range("P456").Value = 8750
$ws.Range("S456").Value = 1250

$ws.Range("D457").Value = 44665
$ws.Range("L457").Value = "Especial"
$ws.Range("M457").Value = 400
$ws.Range("N457").Value = 12500
$ws.Range("O457").Value = 13000
$ws.Range("P457").Value = 12750
$ws.Range("S457").Value = 1821

$ws.Range("D458").Value = 44665
$ws.Range("L458").Value = "Primera"
$ws.Range("M458").Value = 500
$ws.Range("N458").Value = 10500
$ws.Range("O458").Value = 11000
$ws.Range("P458").Value = 10750
$ws.Range("S458").Value = 1536

$ws.Range("D459").Value = 44665
$ws.Range("L459").Value = "Segunda"
$ws.Range("M459").Value = 400
$ws.Range("N459").Value = 8500
$ws.Range("O459").Value = 9000
$ws.Range("P459").Value = 8750
$ws.Range("S459").Value = 1250

$ws.Range("D460").Value = 44510
$ws.Range("L460").Value = "Especial"
$ws.Range("M460").Value = 400
$ws.Range("N460").Value = 12500
$ws.Range("O460").Value = 13000
$ws.Range("P460").Value = 12750
$ws.Range("S460").Value = 1821

$ws.Range("D461").Value = 44510
$ws.Range("L461").Value = "Primera"
$ws.Range("M461").Value = 360
$ws.Range("N461").Value = 10500
$ws.Range("O461").Value = 11000
$ws.Range("P461").Value = 10750
$ws.Range("S461").Value = 1536

$ws.Range("D462").Value = 44510
$ws.Range("L462").Value = "Segunda"
$ws.Range("M462").Value = 300
$ws.Range("N462").Value = 8500
$ws.Range("O462").Value = 9000
$ws.Range("P462").Value = 8750
$ws.Range("S462").Value = 1250

$ws.Range("D463").Value = 44189
$ws.Range("L463").Value = "Especial"
$ws.Range("M463").Value = 300
$ws.Range("N463").Value = 17500
$ws.Range("O463").Value = 18000
$ws.Range("P463").Value = 17750
$ws.Range("S463").Value = 2536

# Append two new rows (464, 465) completing the Coquimbo / Provincia de Melipilla Frutilla block for 44189 (copying fixed columns from row 463, then setting the variable columns)
$ws.Range("A464").Value = $ws.Range("A463").Value()
$ws.Range("B464").Value = $ws.Range("B463").Value()
$ws.Range("C464").Value = $ws.Range("C463").Value()
$ws.Range("D464").Value = 44189
$ws.Range("D464").NumberFormat = $ws.Range("D463").NumberFormat()
$ws.Range("E464").Value = $ws.Range("E463").Value()
$ws.Range("F464").Value = $ws.Range("F463").Value()
$ws.Range("G464").Value = $ws.Range("G463").Value()
$ws.Range("H464").Value = $ws.Range("H463").Value()
$ws.Range("I464").Value = $ws.Range("I463").Value()
$ws.Range("J464").Value = $ws.Range("J463").Value()
$ws.Range("K464").Value = $ws.Range("K463").Value()
$ws.Range("L464").Value = "Primera"
$ws.Range("M464").Value = 300
$ws.Range("N464").Value = 15500
$ws.Range("O464").Value = 16000
$ws.Range("P464").Value = 15750
$ws.Range("Q464").Value = $ws.Range("Q463").Value()
$ws.Range("R464").Value = $ws.Range("R463").Value()
$ws.Range("S464").Value = 2250
$ws.Range("T464").Value = $ws.Range("T463").Value()

$ws.Range("A465").Value = $ws.Range("A463").Value()
$ws.Range("B465").Value = $ws.Range("B463").Value()
$ws.Range("C465").Value = $ws.Range("C463").Value()
$ws.Range("D465").Value = 44189
$ws.Range("D465").NumberFormat = $ws.Range("D463").NumberFormat()
$ws.Range("E465").Value = $ws.Range("E463").Value()
$ws.Range("F465").Value = $ws.Range("F463").Value()
$ws.Range("G465").Value = $ws.Range("G463").Value()
$ws.Range("H465").Value = $ws.Range("H463").Value()
$ws.Range("I465").Value = $ws.Range("I463").Value()
$ws.Range("J465").Value = $ws.Range("J463").Value()
$ws.Range("K465").Value = $ws.Range("K463").Value()
$ws.Range("L465").Value = "Segunda"
$ws.Range("M465").Value = 240
$ws.Range("N465").Value = 12000
$ws.Range("O465").Value = 12500
$ws.Range("P465").Value = 12250
$ws.Range("Q465").Value = $ws.Range("Q463").Value()
$ws.Range("R465").Value = $ws.Range("R463").Value()
$ws.Range("S465").Value = 1750
$ws.Range("T465").Value = $ws.Range("T463").Value()
